# Updated cryptos list on Thu Oct 19 19:58:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to keep its value as plain text (prevents Excel from
    # re-interpreting numeric-looking strings such as "210.21" as numbers),
    # then restore the default "Normal" style so no stray formatting is left
    # behind on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.687.76"
$ws.Range("E2").Value = "  +1.16%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.566.26"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.30%  "

# Row 5 - BNB
Set-TextValue "D5" "210.21"
$ws.Range("E5").Value = "  -0.38%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.57%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.37%  "

# Row 8 - Solana
Set-TextValue "D8" "25.18"
$ws.Range("E8").Value = "  +5.75%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.245"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.02%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.10%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.789.85"
$ws.Range("E12").Value = "  -0.15%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.576.84"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - WrappedBTC
Set-TextValue "D14" "28.673.08"
$ws.Range("E14").Value = "  +1.17%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.86%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -0.72%  "

# Row 17 - Litecoin
Set-TextValue "D17" "61.34"
$ws.Range("E17").Value = "  +0.62%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "229.04"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19 - Chainlink
Set-TextValue "D19" "7.36"
$ws.Range("E19").Value = "  -0.35%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0680"
$ws.Range("E20").Value = "  -0.14%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.52%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +1.04%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.94%  "

# Row 25 - Monero
Set-TextValue "D25" "151.75"
$ws.Range("E25").Value = "  +0.64%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "14.79"
$ws.Range("E26").Value = "  -0.73%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.50%  "

# Row 28 - BinanceUSD
Set-TextValue "D28" "0.998"

# Row 29 - Cosmos
Set-TextValue "D29" "6.24"
$ws.Range("E29").Value = "  -1.51%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.0461"
$ws.Range("E30").Value = "  -4.00%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -2.20%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.17%  "

# Row 33 - Maker
Set-TextValue "D33" "1.394.12"
$ws.Range("E33").Value = "  +1.14%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -3.06%  "

# Row 35 - TrustWalletToken
Set-TextValue "D35" "1.03"
$ws.Range("E35").Value = "  -4.30%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -1.27%  "

# Row 37 - MXToken
$ws.Range("E37").Value = "  +1.50%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  -2.38%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.78%  "

# Row 40 - RenderToken
$ws.Range("E40").Value = "  +1.33%  "

# Row 41 - ImmutableX
Set-TextValue "D41" "0.520"
$ws.Range("E41").Value = "  +0.09%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.34%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  -1.29%  "

# Row 44 - Kaspa
Set-TextValue "D44" "0.0460"
$ws.Range("E44").Value = "  -3.21%  "

# Row 45 - Aave
Set-TextValue "D45" "64.08"
$ws.Range("E45").Value = "  +2.89%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  -1.51%  "

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.701.76"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48 - WEMIXToken
Set-TextValue "D48" "0.870"
$ws.Range("E48").Value = "  -5.00%  "

# Row 49 - Quant
Set-TextValue "D49" "85.20"
$ws.Range("E49").Value = "  -0.20%  "

# Row 50 - BitcoinSV
Set-TextValue "D50" "43.29"
$ws.Range("E50").Value = "  +7.02%  "

# Row 51 - BabyDogeCoin -> Cronos (coin dropped out, replaced by Cronos)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.0513"
$ws.Range("E51").Value = "  -0.28%  "
